# 439-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-B-EarlyRePayment-Loanproduct4
# Rename/populate the loan product with the real product data and add the
# RBI accounting-mapping rows used by the MIFOS strategy test case.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoan_Input")
$ws2 = $wb.Worksheets.Item("ProductLoan_Output")

$productName = "439-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-B-EarlyRePayment"

# --- Update product identity on both sheets -------------------------------
$ws1.Range("B1").Value = $productName
$ws2.Range("B1").Value = $productName

# shortname is now the numeric product id, not the text "kar6"
$ws1.Range("B3").Value = 439

# nominalinterestratedefault corrected
$ws1.Range("B11").Value = 1

# --- Append the new accounting-mapping rows (29-40) ------------------------
$labels = @(
    "fundsource",
    "loanprotfolio",
    "interestreceivable",
    "penaltiesreceivable",
    "transferinsuspense",
    "feesreceivable",
    "incomefrominterest",
    "incomefrompenalties",
    "incomefromfees",
    "incomefromrecoveryrepayments",
    "loseswrittenoff",
    "overpaymentliability"
)

$values = @(
    "Cash",
    "Loan portfolio ",
    "Interest Receivable ",
    "Penalties Receivable ",
    "Transfer in Suspence ",
    "Fees Receivable",
    "Income from interest",
    "Income from penalties",
    "Income from fees",
    "Income from recovery repayments",
    "Losses Writtenoff ",
    "Overpayment Liability"
)

$startRow = 29

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws1.Cells.Item($startRow + $i, 1).Value = $labels[$i]
}
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws1.Cells.Item($startRow + $i, 2).Value = $values[$i]
}

# Match the existing label/value formatting (column A = grey header style,
# column B = green value style) by copying the format from existing rows.
$ws1.Range("A9").Copy() | Out-Null
$ws1.Range("A" + $startRow + ":A" + ($startRow + $labels.Count - 1)).PasteSpecial(-4122) | Out-Null

$ws1.Range("B10").Copy() | Out-Null
$ws1.Range("B" + $startRow + ":B" + ($startRow + $values.Count - 1)).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- View / selection state -------------------------------------------------
# ProductLoan_Input is scrolled down and selected at A26 ...
$ws1.Activate() | Out-Null
$ws1.Range("A26").Select() | Out-Null

# ... while ProductLoan_Output ends up the active (visible) tab, selected at E15
$ws2.Activate() | Out-Null
$ws2.Range("E15").Select() | Out-Null
